# Update TPM-derived NATMI edge statistics for the Tff3-Cxcr4 ligand-receptor pair sheet.
# Columns A-D (Sending cluster / Ligand / Receptor / Target cluster) are unchanged;
# only the numeric statistics in columns E:T (rows 2-13) are refreshed with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.461763666666667"
$ws.Range("H2").Value = [double]"4.385291"
$ws.Range("I2").Value = [double]"0.4829359810344849"
$ws.Range("J2").Value = [double]"0.482935981034485"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"25.69910333333333"
$ws.Range("N2").Value = [double]"77.09730999999999"
$ws.Range("O2").Value = [double]"0.08761243344445813"
$ws.Range("P2").Value = [double]"0.08761243344445814"
$ws.Range("Q2").Value = [double]"37.56601551857889"
$ws.Range("R2").Value = [double]"338.09413966721"
$ws.Range("S2").Value = [double]"0.04231119649631791"
$ws.Range("T2").Value = [double]"0.04231119649631791"

# Row 3
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.461763666666667"
$ws.Range("H3").Value = [double]"4.385291"
$ws.Range("I3").Value = [double]"0.4829359810344849"
$ws.Range("J3").Value = [double]"0.482935981034485"
$ws.Range("K3").Value = [double]"1"
$ws.Range("L3").Value = [double]"0.3333333333333333"
$ws.Range("M3").Value = [double]"0.1622346666666667"
$ws.Range("N3").Value = [double]"0.486704"
$ws.Range("O3").Value = [double]"0.0005530844306649811"
$ws.Range("P3").Value = [double]"0.0005530844306649812"
$ws.Range("Q3").Value = [double]"0.2371487412071112"
$ws.Range("R3").Value = [double]"2.134338670864"
$ws.Range("S3").Value = [double]"0.0002671043721180922"
$ws.Range("T3").Value = [double]"0.0002671043721180923"

# Row 4
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.461763666666667"
$ws.Range("H4").Value = [double]"4.385291"
$ws.Range("I4").Value = [double]"0.4829359810344849"
$ws.Range("J4").Value = [double]"0.482935981034485"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"46.92720933333334"
$ws.Range("N4").Value = [double]"140.781628"
$ws.Range("O4").Value = [double]"0.1599825079935015"
$ws.Range("P4").Value = [double]"0.1599825079935016"
$ws.Range("Q4").Value = [double]"68.59648958152758"
$ws.Range("R4").Value = [double]"617.3684062337481"
$ws.Range("S4").Value = [double]"0.07726130944619898"
$ws.Range("T4").Value = [double]"0.07726130944619901"

# Row 5
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"1.461763666666667"
$ws.Range("H5").Value = [double]"4.385291"
$ws.Range("I5").Value = [double]"0.4829359810344849"
$ws.Range("J5").Value = [double]"0.482935981034485"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"220.538579"
$ws.Range("N5").Value = [double]"661.615737"
$ws.Range("O5").Value = [double]"0.7518519741313753"
$ws.Range("P5").Value = [double]"0.7518519741313754"
$ws.Range("Q5").Value = [double]"322.3752818804964"
$ws.Range("R5").Value = [double]"2901.377536924467"
$ws.Range("S5").Value = [double]"0.3630963707198499"
$ws.Range("T5").Value = [double]"0.36309637071985"

# Row 6
$ws.Range("E6").Value = [double]"2"
$ws.Range("F6").Value = [double]"0.6666666666666666"
$ws.Range("G6").Value = [double]"0.4482143333333333"
$ws.Range("H6").Value = [double]"1.344643"
$ws.Range("I6").Value = [double]"0.1480805917660089"
$ws.Range("J6").Value = [double]"0.1480805917660089"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"25.69910333333333"
$ws.Range("N6").Value = [double]"77.09730999999999"
$ws.Range("O6").Value = [double]"0.08761243344445813"
$ws.Range("P6").Value = [double]"0.08761243344445814"
$ws.Range("Q6").Value = [double]"11.51870646781444"
$ws.Range("R6").Value = [double]"103.66835821033"
$ws.Range("S6").Value = [double]"0.01297370099051543"
$ws.Range("T6").Value = [double]"0.01297370099051543"

# Row 7
$ws.Range("E7").Value = [double]"2"
$ws.Range("F7").Value = [double]"0.6666666666666666"
$ws.Range("G7").Value = [double]"0.4482143333333333"
$ws.Range("H7").Value = [double]"1.344643"
$ws.Range("I7").Value = [double]"0.1480805917660089"
$ws.Range("J7").Value = [double]"0.1480805917660089"
$ws.Range("K7").Value = [double]"1"
$ws.Range("L7").Value = [double]"0.3333333333333333"
$ws.Range("M7").Value = [double]"0.1622346666666667"
$ws.Range("N7").Value = [double]"0.486704"
$ws.Range("O7").Value = [double]"0.0005530844306649811"
$ws.Range("P7").Value = [double]"0.0005530844306649812"
$ws.Range("Q7").Value = [double]"0.07271590296355555"
$ws.Range("R7").Value = [double]"0.654443126672"
$ws.Range("S7").Value = [double]"8.190106978943652E-05"
$ws.Range("T7").Value = [double]"8.190106978943653E-05"

# Row 8
$ws.Range("E8").Value = [double]"2"
$ws.Range("F8").Value = [double]"0.6666666666666666"
$ws.Range("G8").Value = [double]"0.4482143333333333"
$ws.Range("H8").Value = [double]"1.344643"
$ws.Range("I8").Value = [double]"0.1480805917660089"
$ws.Range("J8").Value = [double]"0.1480805917660089"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"46.92720933333334"
$ws.Range("N8").Value = [double]"140.781628"
$ws.Range("O8").Value = [double]"0.1599825079935015"
$ws.Range("P8").Value = [double]"0.1599825079935016"
$ws.Range("Q8").Value = [double]"21.03344784653378"
$ws.Range("R8").Value = [double]"189.301030618804"
$ws.Range("S8").Value = [double]"0.02369030445588795"
$ws.Range("T8").Value = [double]"0.02369030445588795"

# Row 9
$ws.Range("E9").Value = [double]"2"
$ws.Range("F9").Value = [double]"0.6666666666666666"
$ws.Range("G9").Value = [double]"0.4482143333333333"
$ws.Range("H9").Value = [double]"1.344643"
$ws.Range("I9").Value = [double]"0.1480805917660089"
$ws.Range("J9").Value = [double]"0.1480805917660089"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"220.538579"
$ws.Range("N9").Value = [double]"661.615737"
$ws.Range("O9").Value = [double]"0.7518519741313753"
$ws.Range("P9").Value = [double]"0.7518519741313754"
$ws.Range("Q9").Value = [double]"98.84855216076566"
$ws.Range("R9").Value = [double]"889.6369694468909"
$ws.Range("S9").Value = [double]"0.111334685249816"
$ws.Range("T9").Value = [double]"0.1113346852498161"

# Row 10
$ws.Range("E10").Value = [double]"2"
$ws.Range("F10").Value = [double]"0.6666666666666666"
$ws.Range("G10").Value = [double]"1.116849"
$ws.Range("H10").Value = [double]"3.350547"
$ws.Range("I10").Value = [double]"0.3689834271995063"
$ws.Range("J10").Value = [double]"0.3689834271995063"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"25.69910333333333"
$ws.Range("N10").Value = [double]"77.09730999999999"
$ws.Range("O10").Value = [double]"0.08761243344445813"
$ws.Range("P10").Value = [double]"0.08761243344445814"
$ws.Range("Q10").Value = [double]"28.70201785873"
$ws.Range("R10").Value = [double]"258.31816072857"
$ws.Range("S10").Value = [double]"0.0323275359576248"
$ws.Range("T10").Value = [double]"0.03232753595762481"

# Row 11
$ws.Range("E11").Value = [double]"2"
$ws.Range("F11").Value = [double]"0.6666666666666666"
$ws.Range("G11").Value = [double]"1.116849"
$ws.Range("H11").Value = [double]"3.350547"
$ws.Range("I11").Value = [double]"0.3689834271995063"
$ws.Range("J11").Value = [double]"0.3689834271995063"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.1622346666666667"
$ws.Range("N11").Value = [double]"0.486704"
$ws.Range("O11").Value = [double]"0.0005530844306649811"
$ws.Range("P11").Value = [double]"0.0005530844306649812"
$ws.Range("Q11").Value = [double]"0.181191625232"
$ws.Range("R11").Value = [double]"1.630724627088"
$ws.Range("S11").Value = [double]"0.0002040789887574524"
$ws.Range("T11").Value = [double]"0.0002040789887574525"

# Row 12
$ws.Range("E12").Value = [double]"2"
$ws.Range("F12").Value = [double]"0.6666666666666666"
$ws.Range("G12").Value = [double]"1.116849"
$ws.Range("H12").Value = [double]"3.350547"
$ws.Range("I12").Value = [double]"0.3689834271995063"
$ws.Range("J12").Value = [double]"0.3689834271995063"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"46.92720933333334"
$ws.Range("N12").Value = [double]"140.781628"
$ws.Range("O12").Value = [double]"0.1599825079935015"
$ws.Range("P12").Value = [double]"0.1599825079935016"
$ws.Range("Q12").Value = [double]"52.41060681672401"
$ws.Range("R12").Value = [double]"471.695461350516"
$ws.Range("S12").Value = [double]"0.0590308940914146"
$ws.Range("T12").Value = [double]"0.05903089409141461"

# Row 13
$ws.Range("E13").Value = [double]"2"
$ws.Range("F13").Value = [double]"0.6666666666666666"
$ws.Range("G13").Value = [double]"1.116849"
$ws.Range("H13").Value = [double]"3.350547"
$ws.Range("I13").Value = [double]"0.3689834271995063"
$ws.Range("J13").Value = [double]"0.3689834271995063"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"220.538579"
$ws.Range("N13").Value = [double]"661.615737"
$ws.Range("O13").Value = [double]"0.7518519741313753"
$ws.Range("P13").Value = [double]"0.7518519741313754"
$ws.Range("Q13").Value = [double]"246.308291417571"
$ws.Range("R13").Value = [double]"2216.774622758139"
$ws.Range("S13").Value = [double]"0.2774209181617094"
$ws.Range("T13").Value = [double]"0.2774209181617094"
